{"js": "// Replace the 100 \"a+b=c\" / \"a-b=c\" answers in the 20x5 practice table\n// with their new values, cell by cell (row-major order), matching the\n// target diff exactly. Addressing is purely positional (row/col index),\n// not text-search based, since several old/new values repeat across the\n// table (e.g. \"72+16=88\" is both an old value in one cell and a new\n// value written into a different, earlier cell) \u2014 a naive find/replace\n// would clobber the wrong cell.\nconst newValues = [\n  [\"72+16=88\", \"67+7=74\", \"83-64=19\", \"25+4=29\", \"3+51=54\"],\n  [\"54+20=74\", \"17+21=38\", \"91-79=12\", \"68-57=11\", \"39-17=22\"],\n  [\"12+53=65\", \"71-30=41\", \"14+68=82\", \"44+21=65\", \"59-0=59\"],\n  [\"83-41=42\", \"77+9=86\", \"27+52=79\", \"84-50=34\", \"47-21=26\"],\n  [\"9+15=24\", \"55-29=26\", \"12+6=18\", \"40+55=95\", \"13+77=90\"],\n  [\"12+58=70\", \"47+16=63\", \"64-10=54\", \"11+5=16\", \"90-69=21\"],\n  [\"25-1=24\", \"44+27=71\", \"67-35=32\", \"9+83=92\", \"88+8=96\"],\n  [\"64-61=3\", \"48-7=41\", \"43-37=6\", \"43+11=54\", \"19+76=95\"],\n  [\"33+51=84\", \"22+74=96\", \"21+65=86\", \"42-41=1\", \"84-53=31\"],\n  [\"14+84=98\", \"53-44=9\", \"39+35=74\", \"88-13=75\", \"87-77=10\"],\n  [\"39+26=65\", \"37-9=28\", \"54+4=58\", \"96-91=5\", \"23-14=9\"],\n  [\"73-30=43\", \"61+24=85\", \"19+34=53\", \"19-6=13\", \"34+47=81\"],\n  [\"84-13=71\", \"35+15=50\", \"41-27=14\", \"97-25=72\", \"42+35=77\"],\n  [\"3+57=60\", \"54-29=25\", \"75-48=27\", \"80-29=51\", \"10+78=88\"],\n  [\"91-39=52\", \"35+23=58\", \"58-12=46\", \"73-67=6\", \"87-81=6\"],\n  [\"8+72=80\", \"28+16=44\", \"73+4=77\", \"96-56=40\", \"81-12=69\"],\n  [\"50+37=87\", \"52-5=47\", \"46-21=25\", \"16-10=6\", \"42+24=66\"],\n  [\"87-46=41\", \"27+45=72\", \"91-49=42\", \"83-70=13\", \"77-25=52\"],\n  [\"86+10=96\", \"18+29=47\", \"67-41=26\", \"41-18=23\", \"42-5=37\"],\n  [\"56+43=99\", \"26+24=50\", \"54+22=76\", \"16+80=96\", \"1+49=50\"],\n];\n\nconst body = context.document.body;\nconst tables = body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nif (tables.items.length === 0) {\n  throw new Error(\"No table found in document body\");\n}\n\nconst table = tables.items[0];\ntable.load(\"rowCount\");\nawait context.sync();\n\nconst rowCount = Math.min(table.rowCount, newValues.length);\nfor (let r = 0; r < rowCount; r++) {\n  const rowValues = newValues[r];\n  for (let c = 0; c < rowValues.length; c++) {\n    table.getCell(r, c).value = rowValues[c];\n  }\n}\n\nawait context.sync();\n", "ps1": "# Replace the 100 \"a+b=c\" / \"a-b=c\" answers in the 20x5 practice table\n# with their new values, cell by cell (row/column index), matching the\n# target diff exactly. Addressing is purely positional (row/col index),\n# not text-search based, since several old/new values repeat across the\n# table (e.g. \"72+16=88\" is both an old value in one cell and a new\n# value written into a different, earlier cell) \u2014 a naive find/replace\n# (Content.Find.Execute ReplaceAll) would clobber the wrong cell.\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n$newValues = @(\n  @(\"72+16=88\", \"67+7=74\", \"83-64=19\", \"25+4=29\", \"3+51=54\"),\n  @(\"54+20=74\", \"17+21=38\", \"91-79=12\", \"68-57=11\", \"39-17=22\"),\n  @(\"12+53=65\", \"71-30=41\", \"14+68=82\", \"44+21=65\", \"59-0=59\"),\n  @(\"83-41=42\", \"77+9=86\", \"27+52=79\", \"84-50=34\", \"47-21=26\"),\n  @(\"9+15=24\", \"55-29=26\", \"12+6=18\", \"40+55=95\", \"13+77=90\"),\n  @(\"12+58=70\", \"47+16=63\", \"64-10=54\", \"11+5=16\", \"90-69=21\"),\n  @(\"25-1=24\", \"44+27=71\", \"67-35=32\", \"9+83=92\", \"88+8=96\"),\n  @(\"64-61=3\", \"48-7=41\", \"43-37=6\", \"43+11=54\", \"19+76=95\"),\n  @(\"33+51=84\", \"22+74=96\", \"21+65=86\", \"42-41=1\", \"84-53=31\"),\n  @(\"14+84=98\", \"53-44=9\", \"39+35=74\", \"88-13=75\", \"87-77=10\"),\n  @(\"39+26=65\", \"37-9=28\", \"54+4=58\", \"96-91=5\", \"23-14=9\"),\n  @(\"73-30=43\", \"61+24=85\", \"19+34=53\", \"19-6=13\", \"34+47=81\"),\n  @(\"84-13=71\", \"35+15=50\", \"41-27=14\", \"97-25=72\", \"42+35=77\"),\n  @(\"3+57=60\", \"54-29=25\", \"75-48=27\", \"80-29=51\", \"10+78=88\"),\n  @(\"91-39=52\", \"35+23=58\", \"58-12=46\", \"73-67=6\", \"87-81=6\"),\n  @(\"8+72=80\", \"28+16=44\", \"73+4=77\", \"96-56=40\", \"81-12=69\"),\n  @(\"50+37=87\", \"52-5=47\", \"46-21=25\", \"16-10=6\", \"42+24=66\"),\n  @(\"87-46=41\", \"27+45=72\", \"91-49=42\", \"83-70=13\", \"77-25=52\"),\n  @(\"86+10=96\", \"18+29=47\", \"67-41=26\", \"41-18=23\", \"42-5=37\"),\n  @(\"56+43=99\", \"26+24=50\", \"54+22=76\", \"16+80=96\", \"1+49=50\")\n)\n\n$rowCount = [Math]::Min($t.Rows.Count, $newValues.Count)\nfor ($r = 0; $r -lt $rowCount; $r++) {\n  $rowValues = $newValues[$r]\n  for ($c = 0; $c -lt $rowValues.Count; $c++) {\n    $t.Cell($r + 1, $c + 1).Range.Text = $rowValues[$c]\n  }\n}\n"}
